$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last used row in the data table and append a new row after it
$lastRow = $ws.UsedRange.Rows.Count
$newRow = $lastRow + 1

# Copy formatting (date number format / style) from the cell above in column A
$ws.Cells.Item($lastRow, 1).Copy()
$ws.Cells.Item($newRow, 1).PasteSpecial(-4122)

$ws.Cells.Item($newRow, 1).Value = 42624.619733796295
$ws.Cells.Item($newRow, 2).Value = 78
$ws.Cells.Item($newRow, 3).Value = 0
$ws.Cells.Item($newRow, 4).Value = 0
$ws.Cells.Item($newRow, 5).Value = 0
$ws.Cells.Item($newRow, 6).Value = 0
$ws.Cells.Item($newRow, 7).Value = 0
$ws.Cells.Item($newRow, 8).Value = 0
$ws.Cells.Item($newRow, 9).Value = 0
$ws.Cells.Item($newRow, 10).Value = 0
$ws.Cells.Item($newRow, 11).Value = 0
$ws.Cells.Item($newRow, 12).Value = 0
$ws.Cells.Item($newRow, 13).Value = 0
$ws.Cells.Item($newRow, 14).Value = "Random"
